$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows before row 670 (shifts existing rows 670-717 down to 673-720)
$ws.Rows("670:672").Insert()

# Row 670: Pintón, new week (44746)
$ws.Range("A670").Value = 8
$ws.Range("B670").Value = "Terminal La Palmera de La Serena"
$ws.Range("C670").Value = "Coquimbo"
$ws.Range("D670").Value = 44746
$ws.Range("E670").Value = 4
$ws.Range("F670").Value = "Fruta"
$ws.Range("G670").Value = 100108
$ws.Range("H670").Value = "Tropicales y subtropicales"
$ws.Range("I670").Value = 100108006
$ws.Range("J670").Value = "Plátano"
$ws.Range("K670").Value = "Sin especificar"
$ws.Range("L670").Value = "Pintón"
$ws.Range("M670").Value = 80
$ws.Range("N670").Value = 23000
$ws.Range("O670").Value = 23000
$ws.Range("P670").Value = 23000
$ws.Range("Q670").Value = "$/caja 20 kilos"
$ws.Range("R670").Value = "Ecuador"
$ws.Range("S670").Value = 1150
$ws.Range("T670").Value = 20

# Row 671: Primera Maduro, new week (44746)
$ws.Range("A671").Value = 8
$ws.Range("B671").Value = "Terminal La Palmera de La Serena"
$ws.Range("C671").Value = "Coquimbo"
$ws.Range("D671").Value = 44746
$ws.Range("E671").Value = 4
$ws.Range("F671").Value = "Fruta"
$ws.Range("G671").Value = 100108
$ws.Range("H671").Value = "Tropicales y subtropicales"
$ws.Range("I671").Value = 100108006
$ws.Range("J671").Value = "Plátano"
$ws.Range("K671").Value = "Sin especificar"
$ws.Range("L671").Value = "Primera Maduro"
$ws.Range("M671").Value = 120
$ws.Range("N671").Value = 25000
$ws.Range("O671").Value = 25000
$ws.Range("P671").Value = 25000
$ws.Range("Q671").Value = "$/caja 20 kilos"
$ws.Range("R671").Value = "Ecuador"
$ws.Range("S671").Value = 1250
$ws.Range("T671").Value = 20

# Row 672: Primera Pintón, new week (44746)
$ws.Range("A672").Value = 8
$ws.Range("B672").Value = "Terminal La Palmera de La Serena"
$ws.Range("C672").Value = "Coquimbo"
$ws.Range("D672").Value = 44746
$ws.Range("E672").Value = 4
$ws.Range("F672").Value = "Fruta"
$ws.Range("G672").Value = 100108
$ws.Range("H672").Value = "Tropicales y subtropicales"
$ws.Range("I672").Value = 100108006
$ws.Range("J672").Value = "Plátano"
$ws.Range("K672").Value = "Sin especificar"
$ws.Range("L672").Value = "Primera Pintón"
$ws.Range("M672").Value = 120
$ws.Range("N672").Value = 26000
$ws.Range("O672").Value = 26000
$ws.Range("P672").Value = 26000
$ws.Range("Q672").Value = "$/caja 20 kilos"
$ws.Range("R672").Value = "Ecuador"
$ws.Range("S672").Value = 1300
$ws.Range("T672").Value = 20
